# Peru_FX.xlsx update — extend the FX history table through row 313 and
# correct the low/close values on the existing last row (310).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 310 (low / close) ---
$ws.Cells.Item(310, 5).Value2 = 3.7041   # E310 low
$ws.Cells.Item(310, 6).Value2 = 3.7051   # F310 close

# --- Append new rows 311-313 ---
$rows = @(
    @{ Row = 311; Date = 45047.33333333334; Open = 3.7026; High = 3.7503; Low = 3.6386; Close = 3.6559; Volume = 0 },
    @{ Row = 312; Date = 45078.33333333334; Open = 3.6559; High = 3.6947; Low = 3.6008; Close = 3.61;   Volume = 0 },
    @{ Row = 313; Date = 45110.33333333334; Open = 3.6089; High = 3.6488; Low = 3.6069; Close = 3.6264; Volume = 0 }
)

# The date column (A) on existing rows uses a dedicated style (border +
# bold + centered + custom date format). Re-use that exact style on the
# new rows by copying formats from A310, instead of setting NumberFormat
# directly (which would mint a brand-new, differently composed style).
$ws.Cells.Item(310, 1).Copy()

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)  # xlPasteFormats
    $ws.Cells.Item($row, 1).Value2 = $r.Date
    $ws.Cells.Item($row, 2).Value2 = "FX_IDC:USDPEN"
    $ws.Cells.Item($row, 3).Value2 = $r.Open
    $ws.Cells.Item($row, 4).Value2 = $r.High
    $ws.Cells.Item($row, 5).Value2 = $r.Low
    $ws.Cells.Item($row, 6).Value2 = $r.Close
    $ws.Cells.Item($row, 7).Value2 = $r.Volume
}

$excel.CutCopyMode = 0
